# Fix Implementation Approach timeline formatting:
# Remove markdown italic syntax "*(Months X-Y)*" and merge it into the
# preceding bold "Phase N: ..." run so the whole header (including the
# month range) is a single bold run with no asterisks.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)
$shp = $s.Shapes.Item(3)
$tr = $shp.TextFrame.TextRange

$replacements = @(
    @{ Old = "Phase 1: Pilot *(Months 1-2)*"; New = "Phase 1: Pilot (Months 1-2)" },
    @{ Old = "Phase 2: Expansion *(Months 3-4)*"; New = "Phase 2: Expansion (Months 3-4)" },
    @{ Old = "Phase 3: Optimization *(Months 5-6)*"; New = "Phase 3: Optimization (Months 5-6)" }
)

foreach ($rep in $replacements) {
    $fullText = $tr.Text
    $idx = $fullText.IndexOf($rep.Old)
    if ($idx -ge 0) {
        $startPos = $idx + 1
        $len = $rep.Old.Length
        $sub = $tr.Characters($startPos, $len)
        $sub.Text = $rep.New
    }
}
